# Natmi following Dr Hou advice
# Rebuild the LR-pairs data table with the updated (corrected) specificity
# values, including a new "ECs" sending/target cluster category.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Sending cluster",
    "Ligand symbol",
    "Receptor symbol",
    "Target cluster",
    "Ligand-expressing cells",
    "Ligand detection rate",
    "Ligand average expression value",
    "Ligand total expression value",
    "Ligand derived specificity of average expression value",
    "Ligand derived specificity of total expression value",
    "Receptor-expressing cells",
    "Receptor detection rate",
    "Receptor average expression value",
    "Receptor total expression value",
    "Receptor derived specificity of average expression value",
    "Receptor derived specificity of total expression value",
    "Edge average expression weight",
    "Edge total expression weight",
    "Edge average expression derived specificity",
    "Edge total expression derived specificity"
)

for ($c = 0; $c -lt $headers.Count; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$rows = @(
    @("ECs","Lamc2","Itgb4","ECs",2,0.6666666666666666,0.4421816666666667,1.326545,0.06026482003168283,0.06026482003168283,3,1,3.106744,9.320231999999999,0.1418337845295607,0.1418337845295607,1.373745239826667,12.36370715844,0.008547587499086453,0.008547587499086453),
    @("ECs","Lamc2","Itgb4","FAPs",2,0.6666666666666666,0.4421816666666667,1.326545,0.06026482003168283,0.06026482003168283,3,1,1.522686333333333,4.568059,0.0695159837141737,0.06951598371417368,0.673303980683889,6.059735826155,0.004189368247860072,0.004189368247860072),
    @("ECs","Lamc2","Itgb4","sCs",2,0.6666666666666666,0.4421816666666667,1.326545,0.06026482003168283,0.06026482003168283,3,1,17.274688,51.824064,0.7886502317562657,0.7886502317562656,7.638550330986668,68.74695297888,0.04752786428473631,0.0475278642847363),
    @("FAPs","Lamc2","Itgb4","ECs",3,1,6.369908666666666,19.109726,0.8681531333236113,0.8681531333236113,3,1,3.106744,9.320231999999999,0.1418337845295607,0.1418337845295607,19.78967553071466,178.107079776432,0.123133444450484,0.123133444450484),
    @("FAPs","Lamc2","Itgb4","FAPs",3,1,6.369908666666666,19.109726,0.8681531333236113,0.8681531333236113,3,1,1.522686333333333,4.568059,0.0695159837141737,0.06951598371417368,9.699372871314889,87.29435584183399,0.06035051907753303,0.06035051907753302),
    @("FAPs","Lamc2","Itgb4","sCs",3,1,6.369908666666666,19.109726,0.8681531333236113,0.8681531333236113,3,1,17.274688,51.824064,0.7886502317562657,0.7886502317562656,110.0381848051627,990.343663246464,0.6846691697955943,0.6846691697955942),
    @("sCs","Lamc2","Itgb4","ECs",3,1,0.5252196666666666,1.575659,0.07158204664470585,0.07158204664470584,3,1,3.106744,9.320231999999999,0.1418337845295607,0.1418337845295607,1.631723048098666,14.685507432888,0.01015275257999017,0.01015275257999017),
    @("sCs","Lamc2","Itgb4","FAPs",3,1,0.5252196666666666,1.575659,0.07158204664470585,0.07158204664470584,3,1,1.522686333333333,4.568059,0.0695159837141737,0.06951598371417368,0.7997448084312222,7.197703275881,0.004976096388780594,0.004976096388780592),
    @("sCs","Lamc2","Itgb4","sCs",3,1,0.5252196666666666,1.575659,0.07158204664470585,0.07158204664470584,3,1,17.274688,51.824064,0.7886502317562657,0.7886502317562656,9.073005873130667,81.657052858176,0.0564531976759351,0.05645319767593508)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowValues = $rows[$r]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowValues[$c]
    }
}
